$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("valueObject")

# Insert a new column before column I (9) - shifts I..O to J..P.
# CopyOrigin=1 (xlFormatFromRightOrBelow) so the new blank column inherits
# the formatting of the column being pushed right (matches Excel's actual
# "Insert" behaviour when columns on both sides differ).
$ws.Columns("I").Insert(-4161, 1)

# Header row: new column I (merged I38:I39) gets the new "必須" (Required) header.
$ws.Range("I38").Value = "必須"

# Mark the sample data row (row 40) as required: checkmark in H40 (existing
# "必須" flag column for embeddedAnnotation) and in the new I40 cell.
$ws.Range("H40").Value = "○"
$ws.Range("I40").Value = "○"

$ws.Range("H35").Select()
